$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (interest count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 248
$ws1.Range("F5").Value = 3895
$ws1.Range("F6").Value = 28

# Sheet "全部类型" (all types) - same events, rows shifted (F8 instead of F6)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 248
$ws4.Range("F5").Value = 3895
$ws4.Range("F8").Value = 28
